$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.596.97'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.43%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.641.01'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.27%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.25'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.21%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.644'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.51%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("E9").Value = '  -3.02%  '

$ws.Range("E10").Value = '  +0.67%  '

$ws.Range("E11").Value = '  -0.88%  '

$ws.Range("E12").Value = '  +0.08%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '28.76'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.10%  '

$ws.Range("E14").Value = '  -3.92%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.116.14'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.22%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.378.06'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.57%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.641.65'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.87%  '

$ws.Range("E18").Value = '  -2.68%  '

$ws.Range("E19").Value = '  -0.67%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.50'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.71%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '348.30'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.04%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.06%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.29'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.98%  '

$ws.Range("E24").Value = '  +8.30%  '

$ws.Range("E25").Value = '  -0.22%  '

$ws.Range("E26").Value = '  -1.65%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '595.42'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +10.67%  '

$ws.Range("E28").Value = '  +1.22%  '

$ws.Range("E29").Value = '  +1.73%  '

$ws.Range("E30").Value = '  -0.95%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.24%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.09'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.98%  '

$ws.Range("E33").Value = '  +0.08%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.69'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.86%  '

$ws.Range("E35").Value = '  -0.68%  '

$ws.Range("E36").Value = '  -0.90%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '20.08'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.97%  '

$ws.Range("E38").Value = '  -0.03%  '

$ws.Range("E39").Value = '  +1.79%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '154.54'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.35%  '

$ws.Range("E41").Value = '  +0.01%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.43'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.43%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '158.77'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.13%  '

$ws.Range("E44").Value = '  -0.92%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '23.58'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.04%  '

$ws.Range("E46").Value = '  +0.56%  '

$ws.Range("E47").Value = '  +0.31%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0256'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.88%  '

$ws.Range("E49").Value = '  +2.22%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.25'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.70%  '

$ws.Range("E51").Value = '  -5.74%  '
